# Update the "Playa Acapulco" listing (row 87) to its new label, and add a
# new listing for "Av. Perú - LIVE" as row 88, with the same shape of data
# as the surrounding Chile / Viña del Mar rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing Chile / Viña del Mar "Playa Acapulco" entry.
$ws.Cells.Item(87, 3).Value = "Playa Acapulco - LIVE"

# Copy row 87's formatting down onto the new row 88 before filling in the
# new row's values, so borders/fill match the rest of the table.
$ws.Range("A87:F87").Copy() | Out-Null
$ws.Range("A88:F88").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(88, 1).Value = "LIVE, SEA"
$ws.Cells.Item(88, 2).Value = "-33.01656803200213, -71.55930945437895"
$ws.Cells.Item(88, 3).Value = "Av. Perú - LIVE"
$ws.Cells.Item(88, 4).Value = "Viña del Mar"
$ws.Cells.Item(88, 5).Value = "Chile"
$ws.Cells.Item(88, 6).Value = "WAL7kBTXMvM"

# Move the active selection to the next empty row, matching where Excel
# would leave the cursor after adding this entry.
$ws.Range("A89").Select() | Out-Null
